# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new blank column before the
# existing "Late" column (column N), pushing "Late" / "heading" / "Outstanding"
# one column to the right. Then make "Repayment schedule" the active
# sheet/tab (which also clears the previous tab-selection on
# "NewLoanInput"), with the last used selection on that sheet.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before column N ("Late"); existing N/O/P shift to O/P/Q.
$wsSchedule.Columns("N:N").Insert()

# The newly inserted column keeps the same rendered width (11 chars) as its
# neighbour (M), but - unlike the surrounding bestFit columns - carries an
# explicit custom width rather than an autofit one.
$wsSchedule.Columns("N:N").ColumnWidth = 10.17

# "Repayment schedule" becomes the active sheet/tab, with its own selection.
$wsSchedule.Activate()
$wsSchedule.Range("R10").Select()

$wb.Save()
